$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 357). The workbook was refreshed, bumping that date
# from 2023-10-04 (serial 45203) to 2023-10-06 (serial 45205) for all
# rows.
$ws.Range("C2:C357").Value = 45205
